# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country labels caused by the data source re-sorting rows ---
# Row 65 was "Ghana" and Row 66 was "Kirguistan" before the refresh; after
# the refresh Kirguistan's total overtakes Ghana's, so their names swap.
$ws.Range("A65").Value = "Kirguistan"
$ws.Range("A66").Value = "Ghana"

# Row 207 was "Timor Oriental" and Row 208 was "Santa Lucia" before the
# refresh; their names swap as well (values for these two rows stay equal).
$ws.Range("A207").Value = "Santa Lucia"
$ws.Range("A208").Value = "Timor Oriental"

# --- Update the "last refreshed" timestamp message ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Septiembre de 2020 a las 08:01"

# --- Refresh numeric figures for the affected rows ---
# Row 28: Ucrania
$ws.Range("B28").Value = 204932
$ws.Range("C28").Value = 3627
$ws.Range("E28").Value = 112414
$ws.Range("G28").Value = 69
$ws.Range("H28").Value = 4065

# Row 60: Uzbekistan
$ws.Range("B60").Value = 56068
$ws.Range("C60").Value = 292
$ws.Range("D60").Value = 52466
$ws.Range("E60").Value = 3140
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 462

# Row 65: now Kirguistan (new, higher totals)
$ws.Range("B65").Value = 46522
$ws.Range("C65").Value = 167
$ws.Range("D65").Value = 42761
$ws.Range("E65").Value = 2697
$ws.Range("H65").Value = 1064

# Row 66: now Ghana (figures unchanged from before, only label moved)
$ws.Range("B66").Value = 46444
$ws.Range("D66").Value = 45646
$ws.Range("E66").Value = 499
$ws.Range("H66").Value = 299
